$d = $word.ActiveDocument

# Locate the paragraph that contains the "red chick weight" figure caption
# (Figure 2 caption), which is the anchor point for the new paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*This is my red chick weight*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find anchor paragraph containing 'This is my red chick weight'"
}

$anchorPara = $d.Paragraphs.Item($targetIndex)

# Insert a new, empty paragraph right after the anchor paragraph.
$anchorPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newRange = $newPara.Range

$bodyXml = '<w:p><w:pPr><w:pStyle w:val="Textkrper"/></w:pPr><w:r><w:t xml:space="preserve">There is a significant difference in the average weight of chicks who received Diet 1 compared to Diet 3,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>Δ</m:t></m:r><m:r><m:t>M</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>=</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>−</m:t></m:r><m:r><m:t>40.30</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">, 95% CI</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:d><m:dPr><m:begChr m:val="["/><m:endChr m:val="]"/><m:sepChr m:val=""/><m:grow/></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>−</m:t></m:r><m:r><m:t>57.62</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>,</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>−</m:t></m:r><m:r><m:t>22.99</m:t></m:r></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve">,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>t</m:t></m:r><m:d><m:dPr><m:begChr m:val="("/><m:endChr m:val=")"/><m:sepChr m:val=""/><m:grow/></m:dPr><m:e><m:r><m:t>175.92</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>=</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4.59</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>p</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>&lt;</m:t></m:r><m:r><m:t>.001</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">.</w:t></w:r></w:p>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($packageXml)

Write-Output ("Inserted new paragraph at index " + ($targetIndex + 1))
